{"js": "// Append a new work-log row to the end of the first (\"Date / Activity / Time\")\n// table: 10/4/24 | Continued working on ICT risk sections | 3\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst logTable = tables.items[0];\n\nlogTable.addRows(\n  \"End\",\n  1,\n  [[\"10/4/24\", \"Continued working on ICT risk sections\", \"3\"]]\n);\n\nawait context.sync();\n", "ps1": "# Append a new work-log row to the end of the first (\"Date / Activity / Time\")\n# table: 10/4/24 | Continued working on ICT risk sections | 3\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newRow = $t.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"10/4/24\"\n$newRow.Cells.Item(2).Range.Text = \"Continued working on ICT risk sections\"\n$newRow.Cells.Item(3).Range.Text = \"3\"\n"}
